$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1 ("This is a Microsoft word document."):
#    - append two trailing spaces to the existing run
#    - append a parenthetical note in red (split across 3 runs, matching
#      the source document's run boundaries)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.InsertAfter("  ")

$seg1 = "(This is a change " + [char]0x2013 + " Ve"
$seg2 = "rsion for branch alternate"
$seg3 = ")"
$redColor = 192   # RGB(192,0,0) == hex C00000

$rA = $d.Paragraphs.Item(1).Range
$posA = $rA.End - 1
$rA.InsertAfter($seg1)
$runA = $d.Range($posA, $posA + $seg1.Length)
$runA.Font.Color = $redColor

$rB = $d.Paragraphs.Item(1).Range
$posB = $rB.End - 1
$rB.InsertAfter($seg2)
$runB = $d.Range($posB, $posB + $seg2.Length)
$runB.Font.Color = $redColor

$rC = $d.Paragraphs.Item(1).Range
$posC = $rC.End - 1
$rC.InsertAfter($seg3)
$runC = $d.Range($posC, $posC + $seg3.Length)
$runC.Font.Color = $redColor

# ---------------------------------------------------------------------------
# 2) Insert a new, empty, shaded paragraph right after paragraph 2
#    ("It will be treated as a binary file by Git.")
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(3)
$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 16382457

$newRange = $newPara.Range
$newRange.Font.NameAscii = "Calibri"
$newRange.Font.NameFarEast = "Times New Roman"
$newRange.Font.Name = "Calibri"
$newRange.Font.NameBi = "Calibri"
$newRange.Font.Bold = 1
$newRange.Font.BoldBi = 1
$newRange.Font.Color = 2236704   # hex 202122

Write-Output "edit complete"
